$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracking")
$ws.Select()
